$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: "_old" -> "_FV2310", "_new" -> "_FV2404"
$headers = @(
  "Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310",
  "Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310",
  "diff",
  "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
  "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the used range into an Excel Table ("ListObject") covering A1:U65
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U65"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row (pane split after row 1)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
